$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 23 new player rows (309-331) for this week's session
$ws.Cells.Item(309, 1).Value = "Digão"
$ws.Range("C309").Value = 0
$ws.Range("D309").Value = 1
$ws.Range("E309").Value = 6
$ws.Range("F309").Value = 0
$ws.Range("G309").Value = 1
$ws.Range("H309").Value = 0
$ws.Range("I309").Value = 1
$ws.Range("J309").Value = 0
$ws.Range("K309").Value = 0

$ws.Cells.Item(310, 1).Value = "Cabeleira"
$ws.Range("C310").Value = 0
$ws.Range("D310").Value = 1
$ws.Range("E310").Value = 6
$ws.Range("F310").Value = 0
$ws.Range("G310").Value = 1
$ws.Range("H310").Value = 0
$ws.Range("I310").Value = 1
$ws.Range("J310").Value = 0
$ws.Range("K310").Value = 0

$ws.Cells.Item(311, 1).Value = "Jorge"
$ws.Range("C311").Value = 0
$ws.Range("D311").Value = 1
$ws.Range("E311").Value = 6
$ws.Range("F311").Value = 0
$ws.Range("G311").Value = 1
$ws.Range("H311").Value = 0
$ws.Range("I311").Value = 1
$ws.Range("J311").Value = 0
$ws.Range("K311").Value = 0

$ws.Cells.Item(312, 1).Value = "Eder"
$ws.Range("C312").Value = 0
$ws.Range("D312").Value = 1
$ws.Range("E312").Value = 6
$ws.Range("F312").Value = 0
$ws.Range("G312").Value = 1
$ws.Range("H312").Value = 0
$ws.Range("I312").Value = 1
$ws.Range("J312").Value = 0
$ws.Range("K312").Value = 0

$ws.Cells.Item(313, 1).Value = "Eduardo"
$ws.Range("C313").Value = 0
$ws.Range("D313").Value = 1
$ws.Range("E313").Value = 6
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 1
$ws.Range("H313").Value = 0
$ws.Range("I313").Value = 1
$ws.Range("J313").Value = 0
$ws.Range("K313").Value = 0

$ws.Cells.Item(314, 1).Value = "Athos"
$ws.Range("C314").Value = 1
$ws.Range("D314").Value = 2
$ws.Range("E314").Value = 5
$ws.Range("F314").Value = 1
$ws.Range("G314").Value = 1
$ws.Range("H314").Value = 0
$ws.Range("I314").Value = 0
$ws.Range("J314").Value = 0
$ws.Range("K314").Value = 0

$ws.Cells.Item(315, 1).Value = "Ismael"
$ws.Range("C315").Value = 1
$ws.Range("D315").Value = 2
$ws.Range("E315").Value = 5
$ws.Range("F315").Value = 0
$ws.Range("G315").Value = 1
$ws.Range("H315").Value = 0
$ws.Range("I315").Value = 0
$ws.Range("J315").Value = 0
$ws.Range("K315").Value = 0

$ws.Cells.Item(316, 1).Value = "Geovane"
$ws.Range("C316").Value = 1
$ws.Range("D316").Value = 2
$ws.Range("E316").Value = 5
$ws.Range("F316").Value = 3
$ws.Range("G316").Value = 1
$ws.Range("H316").Value = 0
$ws.Range("I316").Value = 0
$ws.Range("J316").Value = 0
$ws.Range("K316").Value = 0

$ws.Cells.Item(317, 1).Value = "Leandrinho"
$ws.Range("C317").Value = 1
$ws.Range("D317").Value = 2
$ws.Range("E317").Value = 5
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 1
$ws.Range("H317").Value = 0
$ws.Range("I317").Value = 0
$ws.Range("J317").Value = 0
$ws.Range("K317").Value = 0

$ws.Cells.Item(318, 1).Value = "Guinha"
$ws.Range("C318").Value = 1
$ws.Range("D318").Value = 2
$ws.Range("E318").Value = 5
$ws.Range("F318").Value = 0
$ws.Range("G318").Value = 1
$ws.Range("H318").Value = 0
$ws.Range("I318").Value = 0
$ws.Range("J318").Value = 0
$ws.Range("K318").Value = 0

$ws.Cells.Item(319, 1).Value = "Leah"
$ws.Range("C319").Value = 12
$ws.Range("D319").Value = 1
$ws.Range("E319").Value = 1
$ws.Range("F319").Value = 3
$ws.Range("G319").Value = 1
$ws.Range("H319").Value = 1
$ws.Range("I319").Value = 0
$ws.Range("J319").Value = 0
$ws.Range("K319").Value = 0

$ws.Cells.Item(320, 1).Value = "Romario"
$ws.Range("C320").Value = 12
$ws.Range("D320").Value = 1
$ws.Range("E320").Value = 1
$ws.Range("F320").Value = 7
$ws.Range("G320").Value = 1
$ws.Range("H320").Value = 1
$ws.Range("I320").Value = 0
$ws.Range("J320").Value = 0
$ws.Range("K320").Value = 0

$ws.Cells.Item(321, 1).Value = "Leandrao"
$ws.Range("C321").Value = 12
$ws.Range("D321").Value = 1
$ws.Range("E321").Value = 1
$ws.Range("F321").Value = 3
$ws.Range("G321").Value = 1
$ws.Range("H321").Value = 1
$ws.Range("I321").Value = 0
$ws.Range("J321").Value = 0
$ws.Range("K321").Value = 0

$ws.Cells.Item(322, 1).Value = "Babão"
$ws.Range("C322").Value = 12
$ws.Range("D322").Value = 1
$ws.Range("E322").Value = 1
$ws.Range("F322").Value = 3
$ws.Range("G322").Value = 1
$ws.Range("H322").Value = 1
$ws.Range("I322").Value = 0
$ws.Range("J322").Value = 0
$ws.Range("K322").Value = 0

$ws.Cells.Item(323, 1).Value = "Adriano"
$ws.Range("C323").Value = 12
$ws.Range("D323").Value = 1
$ws.Range("E323").Value = 1
$ws.Range("F323").Value = 5
$ws.Range("G323").Value = 1
$ws.Range("H323").Value = 1
$ws.Range("I323").Value = 0
$ws.Range("J323").Value = 0
$ws.Range("K323").Value = 0

$ws.Cells.Item(324, 1).Value = "Fabinho"
$ws.Range("C324").Value = 4
$ws.Range("D324").Value = 0
$ws.Range("E324").Value = 5
$ws.Range("F324").Value = 3
$ws.Range("G324").Value = 1
$ws.Range("H324").Value = 0
$ws.Range("I324").Value = 0
$ws.Range("J324").Value = 0
$ws.Range("K324").Value = 0

$ws.Cells.Item(325, 1).Value = "Miqueias"
$ws.Range("C325").Value = 4
$ws.Range("D325").Value = 0
$ws.Range("E325").Value = 5
$ws.Range("F325").Value = 2
$ws.Range("G325").Value = 1
$ws.Range("H325").Value = 0
$ws.Range("I325").Value = 0
$ws.Range("J325").Value = 0
$ws.Range("K325").Value = 0

$ws.Cells.Item(326, 1).Value = "Marcos"
$ws.Range("C326").Value = 4
$ws.Range("D326").Value = 0
$ws.Range("E326").Value = 5
$ws.Range("F326").Value = 2
$ws.Range("G326").Value = 1
$ws.Range("H326").Value = 0
$ws.Range("I326").Value = 0
$ws.Range("J326").Value = 0
$ws.Range("K326").Value = 0

$ws.Cells.Item(327, 1).Value = "Corinthiano"
$ws.Range("C327").Value = 4
$ws.Range("D327").Value = 0
$ws.Range("E327").Value = 5
$ws.Range("F327").Value = 2
$ws.Range("G327").Value = 1
$ws.Range("H327").Value = 0
$ws.Range("I327").Value = 0
$ws.Range("J327").Value = 0
$ws.Range("K327").Value = 0

$ws.Cells.Item(328, 1).Value = "Juscielio"
$ws.Range("C328").Value = 4
$ws.Range("D328").Value = 0
$ws.Range("E328").Value = 5
$ws.Range("F328").Value = 0
$ws.Range("G328").Value = 1
$ws.Range("H328").Value = 0
$ws.Range("I328").Value = 0
$ws.Range("J328").Value = 0
$ws.Range("K328").Value = 0

$ws.Cells.Item(329, 1).Value = "Matheus"
$ws.Range("C329").Value = 4
$ws.Range("D329").Value = 2
$ws.Range("E329").Value = 8
$ws.Range("F329").Value = 0
$ws.Range("G329").Value = 1
$ws.Range("H329").Value = 0
$ws.Range("I329").Value = 0
$ws.Range("J329").Value = 0
$ws.Range("K329").Value = 21

$ws.Cells.Item(330, 1).Value = "Chelin"
$ws.Range("C330").Value = 12
$ws.Range("D330").Value = 2
$ws.Range("E330").Value = 2
$ws.Range("F330").Value = 2
$ws.Range("G330").Value = 1
$ws.Range("H330").Value = 1
$ws.Range("I330").Value = 0
$ws.Range("J330").Value = 1
$ws.Range("K330").Value = 6

$ws.Cells.Item(331, 1).Value = "Iuri"
$ws.Range("C331").Value = 1
$ws.Range("D331").Value = 0
$ws.Range("E331").Value = 6
$ws.Range("F331").Value = 0
$ws.Range("G331").Value = 1
$ws.Range("H331").Value = 0
$ws.Range("I331").Value = 1
$ws.Range("J331").Value = 0
$ws.Range("K331").Value = 12

# Update the frozen-pane scroll position and active selection to match
$ws.Range("C331").Select()
